$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-55 (Chemistry -> ENV bank) in place, column by column
$ws.Cells.Item(2, 1).Value = "MQKJY9"
$ws.Cells.Item(2, 2).Value = "The reason for Sea Level rise is"
$ws.Cells.Item(2, 3).Value = "ENV"
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = "Depletion of ozone layer, Global warming, Smog, Acid Rain"
$ws.Cells.Item(2, 6).Value = "MCQ"
$ws.Cells.Item(2, 7).Value = "Global warming"

$ws.Cells.Item(3, 1).Value = "IMAG7Q"
$ws.Cells.Item(3, 2).Value = "The type of generator used in HAWT which supply power to the grid line is"
$ws.Cells.Item(3, 3).Value = "ENV"
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = "Synchronous Induction Generator, D.C Generator, A.C. Generator, Asynchronous Generator"
$ws.Cells.Item(3, 6).Value = "MCQ"
$ws.Cells.Item(3, 7).Value = "Synchronous Induction Generator"

$ws.Cells.Item(4, 1).Value = "JCHNCH"
$ws.Cells.Item(4, 2).Value = "Atmosphere of big cities is polluted by"
$ws.Cells.Item(4, 3).Value = "ENV"
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = "Household Waste, Radioactive Fallout, Pesticides, Automobile Exhaust"
$ws.Cells.Item(4, 6).Value = "MCQ"
$ws.Cells.Item(4, 7).Value = "Automobile Exhaust"

$ws.Cells.Item(5, 1).Value = "7LVESS"
$ws.Cells.Item(5, 2).Value = "BOD Measures"
$ws.Cells.Item(5, 3).Value = "ENV"
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = "content of bacteria, Content of inorganic matter, Carbonic matter in domestic sewage, OUD"
$ws.Cells.Item(5, 6).Value = "MCQ"
$ws.Cells.Item(5, 7).Value = "Carbonic matter in domestic sewage"

$ws.Cells.Item(6, 1).Value = "5BXZLR"
$ws.Cells.Item(6, 2).Value = "BOD stands for"
$ws.Cells.Item(6, 3).Value = "ENV"
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = "Biochemical Oxygen Dissolution, Biochemical Oxygen Demand, Biochemical Oxidation Demand, Biological Oxygen Demand"
$ws.Cells.Item(6, 6).Value = "MCQ"
$ws.Cells.Item(6, 7).Value = "Biochemical Oxygen Demand"

$ws.Cells.Item(7, 1).Value = "77JGTD"
$ws.Cells.Item(7, 2).Value = "The component of environment made of sea, rivers, lakes, etc is called"
$ws.Cells.Item(7, 3).Value = "ENV"
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = "Atmosphere, Hydrosphere, Biosphere, Lithosphere"
$ws.Cells.Item(7, 6).Value = "MCQ"
$ws.Cells.Item(7, 7).Value = "Hydrosphere"

$ws.Cells.Item(8, 1).Value = "5E3MPN"
$ws.Cells.Item(8, 2).Value = "In which atmospheric layer do most cloud form ?"
$ws.Cells.Item(8, 3).Value = "ENV"
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = "Stratosphere, Atmosphere, Mesosphere, Troposphere"
$ws.Cells.Item(8, 6).Value = "MCQ"
$ws.Cells.Item(8, 7).Value = "Troposphere"

$ws.Cells.Item(9, 1).Value = "DDBG0S"
$ws.Cells.Item(9, 2).Value = "Which type of clouds are responible for thunderstorm"
$ws.Cells.Item(9, 3).Value = "ENV"
$ws.Cells.Item(9, 4).Value = 2
$ws.Cells.Item(9, 5).Value = "Cirrus, Stratus, Cumulonimbus, Nacreous"
$ws.Cells.Item(9, 6).Value = "MCQ"
$ws.Cells.Item(9, 7).Value = "Cumulonimbus"

$ws.Cells.Item(10, 1).Value = "H1O54A"
$ws.Cells.Item(10, 2).Value = "Which is the lowest layer of the atmosphere"
$ws.Cells.Item(10, 3).Value = "ENV"
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 5).Value = "Stratosphere, Mesosphere, Troposphere, Thermosphere"
$ws.Cells.Item(10, 6).Value = "MCQ"
$ws.Cells.Item(10, 7).Value = "Troposphere"

$ws.Cells.Item(11, 1).Value = "5UNR49"
$ws.Cells.Item(11, 2).Value = "The ozone layer is mainly found in which atmospheric layer"
$ws.Cells.Item(11, 3).Value = "ENV"
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = "Troposphere, Stratosphere, Mesosphere, Exosphere"
$ws.Cells.Item(11, 6).Value = "MCQ"
$ws.Cells.Item(11, 7).Value = "Stratosphere"

$ws.Cells.Item(12, 1).Value = "GW5AHZ"
$ws.Cells.Item(12, 2).Value = "In which layer do most airplanes fly ?"
$ws.Cells.Item(12, 3).Value = "ENV"
$ws.Cells.Item(12, 4).Value = 2
$ws.Cells.Item(12, 5).Value = "Stratosphere, Troposphere, Thermosphere, Mesosphere"
$ws.Cells.Item(12, 6).Value = "MCQ"
$ws.Cells.Item(12, 7).Value = "Stratosphere"

$ws.Cells.Item(13, 1).Value = "XYEDII"
$ws.Cells.Item(13, 2).Value = "The hottest layer of the atmosphere is"
$ws.Cells.Item(13, 3).Value = "ENV"
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = "Troposphere, Stratosphere, Thermosphere, Mesosphere"
$ws.Cells.Item(13, 6).Value = "MCQ"
$ws.Cells.Item(13, 7).Value = "Thermosphere"

$ws.Cells.Item(14, 1).Value = "C9CPQJ"
$ws.Cells.Item(14, 2).Value = "Which layer is known as the `"ionosphere`" due to its charged particles ?"
$ws.Cells.Item(14, 3).Value = "ENV"
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(14, 5).Value = "Troposphere, Stratosphere, Thermosphere, Mesosphere"
$ws.Cells.Item(14, 6).Value = "MCQ"
$ws.Cells.Item(14, 7).Value = "Thermosphere"

$ws.Cells.Item(15, 1).Value = "6B454N"
$ws.Cells.Item(15, 2).Value = "The outermost layer of the Earth's atmosphere is called ?"
$ws.Cells.Item(15, 3).Value = "ENV"
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(15, 5).Value = "Troposphere, Mesosphere, Thermosphere, Exosphere"
$ws.Cells.Item(15, 6).Value = "MCQ"
$ws.Cells.Item(15, 7).Value = "Exosphere"

$ws.Cells.Item(16, 1).Value = "7VINIG"
$ws.Cells.Item(16, 2).Value = "In a food chain"
$ws.Cells.Item(16, 3).Value = "ENV"
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 5).Value = "energy flows from tropic level-1 to tropic level-4, energy flows from tropic level-4 to tropic level-1, energy does not flow, energy flows from tropic level-4 to tropic level-3"
$ws.Cells.Item(16, 6).Value = "MCQ"
$ws.Cells.Item(16, 7).Value = "energy flows from tropic level-1 to tropic level-4"

$ws.Cells.Item(17, 1).Value = "XHD9O5"
$ws.Cells.Item(17, 2).Value = "Which fuel used in automobiles produce least pollution ?"
$ws.Cells.Item(17, 3).Value = "ENV"
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = "CNG, LPG, Diesel, Petrol"
$ws.Cells.Item(17, 6).Value = "MCQ"
$ws.Cells.Item(17, 7).Value = "CNG"

$ws.Cells.Item(18, 1).Value = "6MVMEP"
$ws.Cells.Item(18, 2).Value = "Coal is the main contributor of"
$ws.Cells.Item(18, 3).Value = "ENV"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = "CO₂, SO₂, N₂, CO"
$ws.Cells.Item(18, 6).Value = "MCQ"
$ws.Cells.Item(18, 7).Value = "CO₂"

$ws.Cells.Item(19, 1).Value = "6WEZCY"
$ws.Cells.Item(19, 2).Value = "Which gas is mainly responsible for green house effect ?"
$ws.Cells.Item(19, 3).Value = "ENV"
$ws.Cells.Item(19, 4).Value = 2
$ws.Cells.Item(19, 5).Value = "CO₂, CFC, CH₄, N₂O"
$ws.Cells.Item(19, 6).Value = "MCQ"
$ws.Cells.Item(19, 7).Value = "CFC"

$ws.Cells.Item(20, 1).Value = "GJ72X7"
$ws.Cells.Item(20, 2).Value = "Hydrogen fuel cell"
$ws.Cells.Item(20, 3).Value = "ENV"
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = "produce electric power by burning of hydrogen gas, works like a battery, Does not require recharging, All of above"
$ws.Cells.Item(20, 6).Value = "MCQ"
$ws.Cells.Item(20, 7).Value = "works like a battery"

$ws.Cells.Item(21, 1).Value = "OG9JP0"
$ws.Cells.Item(21, 2).Value = "What is typical number of blades in most mordern HAWTs ?"
$ws.Cells.Item(21, 3).Value = "ENV"
$ws.Cells.Item(21, 4).Value = 2
$ws.Cells.Item(21, 5).Value = "1, 2, 3, 6"
$ws.Cells.Item(21, 6).Value = "MCQ"
$ws.Cells.Item(21, 7).Value = "3"

$ws.Cells.Item(22, 1).Value = "LM9WCS"
$ws.Cells.Item(22, 2).Value = "What happens if the wind speed exceeds the turbine's cut-out speed ?"
$ws.Cells.Item(22, 3).Value = "ENV"
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 5).Value = "The turbine speeds up, The turbine stops to prevent damage, The turbine stores extra energy, The turbine reverses its rotation"
$ws.Cells.Item(22, 6).Value = "MCQ"
$ws.Cells.Item(22, 7).Value = "The turbine stops to prevent damage"

$ws.Cells.Item(23, 1).Value = "HVNHFE"
$ws.Cells.Item(23, 2).Value = "What is the main function of the yaw mechanism in a HAWT ?"
$ws.Cells.Item(23, 3).Value = "ENV"
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 5).Value = "To adjust the blade pitch, To rotate the nacelle to face the wind, To slow down the rotor speed, To generate electricity directly"
$ws.Cells.Item(23, 6).Value = "MCQ"
$ws.Cells.Item(23, 7).Value = "To rotate the nacelle to face the wind"

$ws.Cells.Item(24, 1).Value = "S79YZI"
$ws.Cells.Item(24, 2).Value = "Where are the most HAWTs installed ?"
$ws.Cells.Item(24, 3).Value = "ENV"
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = "Underground, Offshore and Onshore, Inside buildings, In urban rooftops"
$ws.Cells.Item(24, 6).Value = "MCQ"
$ws.Cells.Item(24, 7).Value = "Offshore and Onshore"

$ws.Cells.Item(25, 1).Value = "LIWGJX"
$ws.Cells.Item(25, 2).Value = "What is the main advantage of a HAWT over a VAWT (Vertical Axis Wind Turbine) ?"
$ws.Cells.Item(25, 3).Value = "ENV"
$ws.Cells.Item(25, 4).Value = 2
$ws.Cells.Item(25, 5).Value = "Can operate in all wind direction without orientation, More efficient at converting wind energy, Require less land area, Works better in low wind speed"
$ws.Cells.Item(25, 6).Value = "MCQ"
$ws.Cells.Item(25, 7).Value = "More efficient at converting wind energy"

$ws.Cells.Item(26, 1).Value = "SGYA00"
$ws.Cells.Item(26, 2).Value = "In horizontal axis Wind Turbine (HAWT) where generator is placed ?"
$ws.Cells.Item(26, 3).Value = "ENV"
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = "In the nessel, Above the tower, At ground level, In the gearbox"
$ws.Cells.Item(26, 6).Value = "MCQ"
$ws.Cells.Item(26, 7).Value = "In the nessel"

$ws.Cells.Item(27, 1).Value = "73Y3YS"
$ws.Cells.Item(27, 2).Value = "Which of the following solid waste disposal methods is ecologically most acceptable ?"
$ws.Cells.Item(27, 3).Value = "ENV"
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = "Sanitary land fill, Inclineration, Composting, Pyrolysis"
$ws.Cells.Item(27, 6).Value = "MCQ"
$ws.Cells.Item(27, 7).Value = "Composting"

$ws.Cells.Item(28, 1).Value = "0T93V5"
$ws.Cells.Item(28, 2).Value = "In a HAWT, the rotor shaft is aligned in which direction ?"
$ws.Cells.Item(28, 3).Value = "ENV"
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = "Horizontal, Vertical, Diagonal, Perpendicular to wind Flow"
$ws.Cells.Item(28, 6).Value = "MCQ"
$ws.Cells.Item(28, 7).Value = "Horizontal"

$ws.Cells.Item(29, 1).Value = "ATY726"
$ws.Cells.Item(29, 2).Value = "What does HAWT stand for ?"
$ws.Cells.Item(29, 3).Value = "ENV"
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).Value = "Horizontal Axis Wind Turbine, Horizontal Axis Wind Tubeline, Horizonal Airflow wind Turbine, Hybrid Axis Wind Technology"
$ws.Cells.Item(29, 6).Value = "MCQ"
$ws.Cells.Item(29, 7).Value = "Horizontal Axis Wind Turbine"

$ws.Cells.Item(30, 1).Value = "B5U4YS"
$ws.Cells.Item(30, 2).Value = "Which is the most common material used for mordern HAWT blades ?"
$ws.Cells.Item(30, 3).Value = "ENV"
$ws.Cells.Item(30, 4).Value = 2
$ws.Cells.Item(30, 5).Value = "Carbon Fibre, Fibreglass - reinforced plastic, Pure aluminium, Titanium Alloy"
$ws.Cells.Item(30, 6).Value = "MCQ"
$ws.Cells.Item(30, 7).Value = "Fibreglass - reinforced plastic"

$ws.Cells.Item(31, 1).Value = "72TN4D"
$ws.Cells.Item(31, 2).Value = "The biggest pollutant receptor or sink on the earth is"
$ws.Cells.Item(31, 3).Value = "ENV"
$ws.Cells.Item(31, 4).Value = 2
$ws.Cells.Item(31, 5).Value = "Hyrdosphere, Lithosphere, Atmosphere, Biosphere"
$ws.Cells.Item(31, 6).Value = "MCQ"
$ws.Cells.Item(31, 7).Value = "Hyrdosphere"

$ws.Cells.Item(32, 1).Value = "9CRL27"
$ws.Cells.Item(32, 2).Value = "As per IS Code, the acceptable noise level in urban residental area is"
$ws.Cells.Item(32, 3).Value = "ENV"
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = "35-45 dBA, 45-55 dBA, 30-40 dBA, 40-50 dBA"
$ws.Cells.Item(32, 6).Value = "MCQ"
$ws.Cells.Item(32, 7).Value = "35-45 dBA"

$ws.Cells.Item(33, 1).Value = "Q3A4ZG"
$ws.Cells.Item(33, 2).Value = "Which is the largest ecosystem on the earth ?"
$ws.Cells.Item(33, 3).Value = "ENV"
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 5).Value = "Forests, Sea, Deserts, Grasslands"
$ws.Cells.Item(33, 6).Value = "MCQ"
$ws.Cells.Item(33, 7).Value = "Sea"

$ws.Cells.Item(34, 1).Value = "IUQBES"
$ws.Cells.Item(34, 2).Value = "In biogas plant, digestion takes place in the absense of"
$ws.Cells.Item(34, 3).Value = "ENV"
$ws.Cells.Item(34, 4).Value = 2
$ws.Cells.Item(34, 5).Value = "Oxygen, Carbon dioxide, Hydrogen, Methane"
$ws.Cells.Item(34, 6).Value = "MCQ"
$ws.Cells.Item(34, 7).Value = "Oxygen"

$ws.Cells.Item(35, 1).Value = "SGY0BX"
$ws.Cells.Item(35, 2).Value = "The daily cover of MSW landfills consists of which of the following"
$ws.Cells.Item(35, 3).Value = "ENV"
$ws.Cells.Item(35, 4).Value = 2
$ws.Cells.Item(35, 5).Value = "Compacted Soil, Geomembrane, Geotextile, Geocomposite"
$ws.Cells.Item(35, 6).Value = "MCQ"
$ws.Cells.Item(35, 7).Value = "Compacted Soil"

$ws.Cells.Item(36, 1).Value = "ADHUWA"
$ws.Cells.Item(36, 2).Value = "In the atmosphere reaction between Nitrogen oxides, hydrocarbon and sunlight produces"
$ws.Cells.Item(36, 3).Value = "ENV"
$ws.Cells.Item(36, 4).Value = 2
$ws.Cells.Item(36, 5).Value = "Sulphur Dioxide, Benzene, CO, PAN"
$ws.Cells.Item(36, 6).Value = "MCQ"
$ws.Cells.Item(36, 7).Value = "PAN"

$ws.Cells.Item(37, 1).Value = "WUAHUD"
$ws.Cells.Item(37, 2).Value = "Noise pollution is measured in"
$ws.Cells.Item(37, 3).Value = "ENV"
$ws.Cells.Item(37, 4).Value = 2
$ws.Cells.Item(37, 5).Value = "Decibel, Newton, Dyne, ECO"
$ws.Cells.Item(37, 6).Value = "MCQ"
$ws.Cells.Item(37, 7).Value = "Decibel"

$ws.Cells.Item(38, 1).Value = "Z7NVHM"
$ws.Cells.Item(38, 2).Value = "Vegetables, grass, trees, etc are which type of biotic components ?"
$ws.Cells.Item(38, 3).Value = "ENV"
$ws.Cells.Item(38, 4).Value = 2
$ws.Cells.Item(38, 5).Value = "Producers, Consumers, Decomposers, Transformers"
$ws.Cells.Item(38, 6).Value = "MCQ"
$ws.Cells.Item(38, 7).Value = "Producers"

$ws.Cells.Item(39, 1).Value = "R8YXIP"
$ws.Cells.Item(39, 2).Value = "Nacelle is used in"
$ws.Cells.Item(39, 3).Value = "ENV"
$ws.Cells.Item(39, 4).Value = 2
$ws.Cells.Item(39, 5).Value = "Darious wind turbine, Savonious wind turbine, Horizontal axis wind turbine, Water pumping wind mill"
$ws.Cells.Item(39, 6).Value = "MCQ"
$ws.Cells.Item(39, 7).Value = "Horizontal axis wind turbine"

$ws.Cells.Item(40, 1).Value = "932EXG"
$ws.Cells.Item(40, 2).Value = "Due to 'Ozone hole' in the earth's atmosphere, which ray comes to earth surface ?"
$ws.Cells.Item(40, 3).Value = "ENV"
$ws.Cells.Item(40, 4).Value = 2
$ws.Cells.Item(40, 5).Value = "Laser rays, X-rays, Ultravoilet rays, None of these"
$ws.Cells.Item(40, 6).Value = "MCQ"
$ws.Cells.Item(40, 7).Value = "Ultravoilet rays"

$ws.Cells.Item(41, 1).Value = "MP134F"
$ws.Cells.Item(41, 2).Value = "Less of removal of the superficial layer of soil by the action of water, wind, or by the human activites are termed as"
$ws.Cells.Item(41, 3).Value = "ENV"
$ws.Cells.Item(41, 4).Value = 2
$ws.Cells.Item(41, 5).Value = "Soil Erosion, Soil Pollution, Desertification, Salination"
$ws.Cells.Item(41, 6).Value = "MCQ"
$ws.Cells.Item(41, 7).Value = "Soil Erosion"

$ws.Cells.Item(42, 1).Value = "OYAT71"
$ws.Cells.Item(42, 2).Value = "Minamata disease is caused by ____"
$ws.Cells.Item(42, 3).Value = "ENV"
$ws.Cells.Item(42, 4).Value = 2
$ws.Cells.Item(42, 5).Value = "Mercury, Lead, Zinc, Manganese"
$ws.Cells.Item(42, 6).Value = "MCQ"
$ws.Cells.Item(42, 7).Value = "Mercury"

$ws.Cells.Item(43, 1).Value = "27WJIU"
$ws.Cells.Item(43, 2).Value = "Which of the following is an artifical ecosystem ?"
$ws.Cells.Item(43, 3).Value = "ENV"
$ws.Cells.Item(43, 4).Value = 2
$ws.Cells.Item(43, 5).Value = "Forests, Desert, Fish house, Pond"
$ws.Cells.Item(43, 6).Value = "MCQ"
$ws.Cells.Item(43, 7).Value = "Fish house"

$ws.Cells.Item(44, 1).Value = "XVPYOI"
$ws.Cells.Item(44, 2).Value = "Which pollutant is primarily responsible for causing respiratory problems and cardiovascular diseases in humans ?"
$ws.Cells.Item(44, 3).Value = "ENV"
$ws.Cells.Item(44, 4).Value = 2
$ws.Cells.Item(44, 5).Value = "Lead, Benzene, Carbon Monoxide (CO), Particulate matter (PM10)"
$ws.Cells.Item(44, 6).Value = "MCQ"
$ws.Cells.Item(44, 7).Value = "Particulate matter (PM10)"

$ws.Cells.Item(45, 1).Value = "WH53C3"
$ws.Cells.Item(45, 2).Value = "What is the approximately range of mordern solar photovoltaic (PV) panels in converting sunlight in to electricity ?"
$ws.Cells.Item(45, 3).Value = "ENV"
$ws.Cells.Item(45, 4).Value = 2
$ws.Cells.Item(45, 5).Value = "5-10%, 15-20%, 20-30%, 40-50%"
$ws.Cells.Item(45, 6).Value = "MCQ"
$ws.Cells.Item(45, 7).Value = "15-20%"

$ws.Cells.Item(46, 1).Value = "1MI8GU"
$ws.Cells.Item(46, 2).Value = "CFCs are not recommended to be used in refridgerator because they ____"
$ws.Cells.Item(46, 3).Value = "ENV"
$ws.Cells.Item(46, 4).Value = 2
$ws.Cells.Item(46, 5).Value = "Increase Temperature, Deplete Ozone, Affect Environment, Affect Human Body"
$ws.Cells.Item(46, 6).Value = "MCQ"
$ws.Cells.Item(46, 7).Value = "Deplete Ozone"

$ws.Cells.Item(47, 1).Value = "F4FMI3"
$ws.Cells.Item(47, 2).Value = "Clouds are present ____ in layer of atmosphere."
$ws.Cells.Item(47, 3).Value = "ENV"
$ws.Cells.Item(47, 4).Value = 2
$ws.Cells.Item(47, 5).Value = "Stratosphere, Troposphere, Thermosphere, Mesosphere"
$ws.Cells.Item(47, 6).Value = "MCQ"
$ws.Cells.Item(47, 7).Value = "Troposphere"

$ws.Cells.Item(48, 1).Value = "KKDLYR"
$ws.Cells.Item(48, 2).Value = "Which pollutant is commonly associated with blue baby syndrome in infants ?"
$ws.Cells.Item(48, 3).Value = "ENV"
$ws.Cells.Item(48, 4).Value = 2
$ws.Cells.Item(48, 5).Value = "Lead, Arsenic, Mercury, Nitrate"
$ws.Cells.Item(48, 6).Value = "MCQ"
$ws.Cells.Item(48, 7).Value = "Nitrate"

$ws.Cells.Item(49, 1).Value = "UWO0O6"
$ws.Cells.Item(49, 2).Value = "Which disease is known as 'pain-pain-diseases' due to serve bone pain and fractures ?"
$ws.Cells.Item(49, 3).Value = "ENV"
$ws.Cells.Item(49, 4).Value = 2
$ws.Cells.Item(49, 5).Value = "Black Lung Disease, Itai-Itai Disease, Sillicosis, Fluorosis"
$ws.Cells.Item(49, 6).Value = "MCQ"
$ws.Cells.Item(49, 7).Value = "Itai-Itai Disease"

$ws.Cells.Item(50, 1).Value = "C2T882"
$ws.Cells.Item(50, 2).Value = "Which occupational diseases is caused by inhaling coal dust over a long period ?"
$ws.Cells.Item(50, 3).Value = "ENV"
$ws.Cells.Item(50, 4).Value = 2
$ws.Cells.Item(50, 5).Value = "Asbestosis, Silicosis, Black Lung Diseases, Eutrophication"
$ws.Cells.Item(50, 6).Value = "MCQ"
$ws.Cells.Item(50, 7).Value = "Black Lung Diseases"

$ws.Cells.Item(51, 1).Value = "65P9UL"
$ws.Cells.Item(51, 2).Value = "What is the primary poltant responsivle for Flurosis ?"
$ws.Cells.Item(51, 3).Value = "ENV"
$ws.Cells.Item(51, 4).Value = 2
$ws.Cells.Item(51, 5).Value = "Fluroide, Lead, Sulfur dioxide, Arsenic"
$ws.Cells.Item(51, 6).Value = "MCQ"
$ws.Cells.Item(51, 7).Value = "Fluroide"

$ws.Cells.Item(52, 1).Value = "IAI21L"
$ws.Cells.Item(52, 2).Value = "The BHOPAL GAS TRAGEDY was caused due to which leakage of gas ?"
$ws.Cells.Item(52, 3).Value = "ENV"
$ws.Cells.Item(52, 4).Value = 2
$ws.Cells.Item(52, 5).Value = "Carbon Monoxide, Methyl isocyanate (MIC), Sulfur dioxide, Nitrogen dioxide"
$ws.Cells.Item(52, 6).Value = "MCQ"
$ws.Cells.Item(52, 7).Value = "Methyl isocyanate (MIC)"

$ws.Cells.Item(53, 1).Value = "481PMN"
$ws.Cells.Item(53, 2).Value = "Which disease is associated with inhaling silica dust ?"
$ws.Cells.Item(53, 3).Value = "ENV"
$ws.Cells.Item(53, 4).Value = 2
$ws.Cells.Item(53, 5).Value = "Silicosis, Asbestosis, Itai-Itai Disease, Blue Baby Syndrome"
$ws.Cells.Item(53, 6).Value = "MCQ"
$ws.Cells.Item(53, 7).Value = "Silicosis"

$ws.Cells.Item(54, 1).Value = "S2UOWF"
$ws.Cells.Item(54, 2).Value = "Acid Rain is mainly caused by which pollutants ?"
$ws.Cells.Item(54, 3).Value = "ENV"
$ws.Cells.Item(54, 4).Value = 2
$ws.Cells.Item(54, 5).Value = "Carbon monoxide and Lead, Nitrogen oxide and Sulphur dioxide, Fluoride and Mercury, Ozone and PM2.5"
$ws.Cells.Item(54, 6).Value = "MCQ"
$ws.Cells.Item(54, 7).Value = "Nitrogen oxide and Sulphur dioxide"

$ws.Cells.Item(55, 1).Value = "ALJDJ9"
$ws.Cells.Item(55, 2).Value = "What is the main source of Cadmium Contamination leading to Itai-Itai diseases"
$ws.Cells.Item(55, 3).Value = "ENV"
$ws.Cells.Item(55, 4).Value = 2
$ws.Cells.Item(55, 5).Value = "Fertilizer runoff, Industrial wastewater, burning of fossil fuels, Oil spills"
$ws.Cells.Item(55, 6).Value = "MCQ"
$ws.Cells.Item(55, 7).Value = "Industrial wastewater"

# Remove now-unused trailing rows 56-70 (old Chemistry bank had more rows than new ENV bank)
$ws.Range("A56:G70").EntireRow.Delete()

